# The workbook was simply re-saved by Excel; cell data, shared strings and
# the table/column definitions are unchanged. The only real content changes
# are that column B got wider and the frozen header row / active selection
# need to be (re)written so the sheetView round-trips the way Excel does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B widened (was 37.63 chars, now ~54.33 chars) ---------------
$ws.Columns.Item(2).ColumnWidth = 53.5

# --- Re-establish the frozen header row + final selection ---------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B3").Select() | Out-Null
